$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 7: Date, Activity, Hours
# Copy the date format from A6 so the new date cell uses the same style
$ws.Range("A6").Copy()
$ws.Range("A7").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

$ws.Range("A7").Value = 44903
$ws.Range("B7").Value = "Trying various approaches for spatial ACF"
$ws.Range("C7").Value = 3

# Update selection to C2:C7 with active cell C2
$ws.Range("C2:C7").Select()
